$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.173.15"
$ws.Range("E2").Value = "  -4.30%  "
$ws.Range("D3").Value = "1.656.99"
$ws.Range("E3").Value = "  -2.99%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'216.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.46%  "
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("D7").Value = "'1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'0.2583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").Value = "'0.06423"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("D10").Value = "'19.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.88%  "
$ws.Range("D11").Value = "'0.07768"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "1.661.41"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("D13").Value = "'4.297"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.85%  "
$ws.Range("D14").Value = "1.884.79"
$ws.Range("E14").Value = "  -3.00%  "
$ws.Range("D15").Value = "'0.5537"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.92%  "
$ws.Range("D16").Value = "0.0₅8032"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "'64.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.05%  "
$ws.Range("D18").Value = "26.199.43"
$ws.Range("E18").Value = "  -4.17%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'210.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Value = "'4.410"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.76%  "
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").Value = "'6.012"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'144.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "'1.734"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").Value = "'6.973"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("D30").Value = "'0.05118"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("D31").Value = "'1.249"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("D33").Value = "'3.229"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.82%  "
$ws.Range("D34").Value = "'1.567"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").Value = "'2.753"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.372"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9291"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.5707"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.165.66"
$ws.Range("E39").Value = "  +11.32%  "
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "'0.8396"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "'5.656"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").Value = "'100.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "1.795.27"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("D47").Value = "'0.4541"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "'7.856"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("D51").Value = "'0.05060"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.35%  "
